# custom_distribution_reports.xlsx edit
# - add a new "initial" sheet between "settings" and "model" that tells the
#   survey engine to skip straight to the finalize screen (no init clause),
#   with translations/comments for the delivery (ctp) flow.
# - the new sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook

$settingsSheet = $wb.Worksheets.Item("settings")

# Insert a new worksheet right after "settings" (i.e. before "model").
$initialSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $settingsSheet)
$initialSheet.Name = "initial"

# Header row.
$initialSheet.Range("A1").Value = "clause"
$initialSheet.Range("B1").Value = "type"
$initialSheet.Range("C1").Value = "display.text"
$initialSheet.Range("D1").Value = "comments"

# "do section survey" clause.
$initialSheet.Range("A2").Value = "do section survey"

# "goto _finalize" clause, skipping the finalize screen.
$initialSheet.Range("A3").Value = "goto _finalize"
$initialSheet.Range("D3").Value = "skips the finalize screen where the user chooses to save as incomplete or finalized and instead saves as finalized"

# Wrap text for the table cells, matching each row's populated columns.
$initialSheet.Range("A1:D1").WrapText = $true
$initialSheet.Range("A2:C2").WrapText = $true
$initialSheet.Range("A3:D3").WrapText = $true

# Taller rows so the wrapped comment text is visible.
$initialSheet.Rows.Item(2).RowHeight = 32
$initialSheet.Rows.Item(3).RowHeight = 208

$initialSheet.Range("G3").Select() | Out-Null
